$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 430-431, pushing the existing rows
# (old 430-449) down to 432-451.
$ws.Range("A430:A431").EntireRow.Insert()

# New row 430: weekly entry for 2021-11-09 (serial 44509), "Primera" quality
$ws.Range("A430").Value = 9
$ws.Range("B430").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C430").Value = "Metropolitana"
$ws.Range("D430").Value = 44509
$ws.Range("E430").Value = 13
$ws.Range("F430").Value = 100112023
$ws.Range("G430").Value = "Brócoli"
$ws.Range("H430").Value = "Sin especificar"
$ws.Range("I430").Value = "Primera"
$ws.Range("J430").Value = 6100
$ws.Range("K430").Value = 600
$ws.Range("L430").Value = 700
$ws.Range("M430").Value = 650
$ws.Range("N430").Value = "`$/unidad"
$ws.Range("O430").Value = "Región Metropolitana"
$ws.Range("P430").Value = 650
$ws.Range("Q430").Value = 1
$ws.Range("R430").Value = "Hortaliza"

# New row 431: weekly entry for 2021-11-09 (serial 44509), "Segunda" quality
$ws.Range("A431").Value = 9
$ws.Range("B431").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C431").Value = "Metropolitana"
$ws.Range("D431").Value = 44509
$ws.Range("E431").Value = 13
$ws.Range("F431").Value = 100112023
$ws.Range("G431").Value = "Brócoli"
$ws.Range("H431").Value = "Sin especificar"
$ws.Range("I431").Value = "Segunda"
$ws.Range("J431").Value = 2500
$ws.Range("K431").Value = 500
$ws.Range("L431").Value = 500
$ws.Range("M431").Value = 500
$ws.Range("N431").Value = "`$/unidad"
$ws.Range("O431").Value = "Región Metropolitana"
$ws.Range("P431").Value = 500
$ws.Range("Q431").Value = 1
$ws.Range("R431").Value = "Hortaliza"
